$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A276").Value = 275
$ws.Range("B276").Value = "Saturday, Jan 14"
$ws.Range("C276").Value = "6:15 PM"
$ws.Range("D276").Value = "W61016"
$ws.Range("E276").Value = "Leeds"
$ws.Range("F276").Value = "(LBA)"
$ws.Range("G276").Value = "Wizz Air "
$ws.Range("H276").Value = "A21N"
$ws.Range("I276").Value = "(HA-LZF)"
$ws.Range("J276").Value = "6:05 PM"
$ws.Range("L276").Value = "0 hours, -10 minutes"

$ws.Range("A277").Value = 276
$ws.Range("B277").Value = "Saturday, Jan 14"
$ws.Range("C277").Value = "6:20 PM"
$ws.Range("D277").Value = "W61220"
$ws.Range("E277").Value = "Bergen"
$ws.Range("F277").Value = "(BGO)"
$ws.Range("G277").Value = "Wizz Air "
$ws.Range("H277").Value = "A320"
$ws.Range("I277").Value = "(HA-LWP)"
$ws.Range("J277").Value = "5:53 PM"
$ws.Range("L277").Value = "0 hours, -27 minutes"

$ws.Range("A278").Value = 277
$ws.Range("B278").Value = "Saturday, Jan 14"
$ws.Range("C278").Value = "6:35 PM"
$ws.Range("D278").Value = "FR6390"
$ws.Range("E278").Value = "Dortmund"
$ws.Range("F278").Value = "(DTM)"
$ws.Range("G278").Value = "Ryanair "
$ws.Range("H278").Value = "B738"
$ws.Range("I278").Value = "(SP-RSN)"
$ws.Range("J278").Value = "6:26 PM"
$ws.Range("L278").Value = "0 hours, -9 minutes"

$ws.Range("A279").Value = 278
$ws.Range("B279").Value = "Saturday, Jan 14"
$ws.Range("C279").Value = "7:20 PM"
$ws.Range("D279").Value = "3Z7335"
$ws.Range("E279").Value = "Lanzarote"
$ws.Range("F279").Value = "(ACE)"
$ws.Range("G279").Value = "Smartwings "
$ws.Range("H279").Value = "B38M"
$ws.Range("I279").Value = "(OK-SWC)"
$ws.Range("J279").Value = "7:08 PM"
$ws.Range("L279").Value = "0 hours, -12 minutes"

$ws.Range("A280").Value = 279
$ws.Range("B280").Value = "Saturday, Jan 14"
$ws.Range("C280").Value = "8:10 PM"
$ws.Range("D280").Value = "E44114"
$ws.Range("E280").Value = "Antalya"
$ws.Range("F280").Value = "(AYT)"
$ws.Range("G280").Value = "Enter Air "
$ws.Range("H280").Value = "B738"
$ws.Range("I280").Value = "(SP-ENQ)"
$ws.Range("J280").Value = "7:51 PM"
$ws.Range("L280").Value = "0 hours, -19 minutes"

$ws.Range("A281").Value = 280
$ws.Range("B281").Value = "Saturday, Jan 14"
$ws.Range("C281").Value = "8:35 PM"
$ws.Range("D281").Value = "LPR42"
$ws.Range("E281").Value = "Warsaw"
$ws.Range("F281").Value = "(WAW)"
$ws.Range("G281").Value = "Polish Medical Air Rescue "
$ws.Range("H281").Value = "LJ75"
$ws.Range("I281").Value = "(SP-MXS)"
$ws.Range("J281").Value = "8:22 PM"
$ws.Range("L281").Value = "0 hours, -13 minutes"

$ws.Range("A282").Value = 281
$ws.Range("B282").Value = "Saturday, Jan 14"
$ws.Range("C282").Value = "9:10 PM"
$ws.Range("D282").Value = "FR6392"
$ws.Range("E282").Value = "London"
$ws.Range("F282").Value = "(STN)"
$ws.Range("G282").Value = "Ryanair "
$ws.Range("H282").Value = "B738"
$ws.Range("I282").Value = "(EI-GSH)"
$ws.Range("J282").Value = "9:16 PM"
$ws.Range("L282").Value = "0 hours, 6 minutes"

$ws.Range("A283").Value = 282
$ws.Range("B283").Value = "Saturday, Jan 14"
$ws.Range("C283").Value = "9:47 PM"
$ws.Range("D283").Value = "UNKNOWN"
$ws.Range("E283").Value = "Rovaniemi"
$ws.Range("F283").Value = "(RVN)"
$ws.Range("G283").Value = "Enter Air "
$ws.Range("H283").Value = "B738"
$ws.Range("I283").Value = "(SP-ESC)"
$ws.Range("J283").Value = "9:58 PM"
$ws.Range("L283").Value = "0 hours, 11 minutes"

$ws.Range("A284").Value = 283
$ws.Range("B284").Value = "Saturday, Jan 14"
$ws.Range("C284").Value = "9:55 PM"
$ws.Range("D284").Value = "FR6404"
$ws.Range("E284").Value = "Manchester"
$ws.Range("F284").Value = "(MAN)"
$ws.Range("G284").Value = "Ryanair "
$ws.Range("H284").Value = "B738"
$ws.Range("I284").Value = "(SP-RSB)"
$ws.Range("J284").Value = "10:10 PM"
$ws.Range("L284").Value = "0 hours, 15 minutes"

$ws.Range("A285").Value = 284
$ws.Range("B285").Value = "Saturday, Jan 14"
$ws.Range("C285").Value = "10:20 PM"
$ws.Range("D285").Value = "W61176"
$ws.Range("E285").Value = "Barcelona"
$ws.Range("F285").Value = "(BCN)"
$ws.Range("G285").Value = "Wizz Air "
$ws.Range("H285").Value = "A321"
$ws.Range("I285").Value = "(HA-LTC)"
$ws.Range("J285").Value = "9:56 PM"
$ws.Range("L285").Value = "0 hours, -24 minutes"

$ws.Range("A286").Value = 285
$ws.Range("B286").Value = "Saturday, Jan 14"
$ws.Range("C286").Value = "11:23 PM"
$ws.Range("D286").Value = "E44004"
$ws.Range("E286").Value = "Sharm el-Sheikh"
$ws.Range("F286").Value = "(SSH)"
$ws.Range("G286").Value = "Enter Air "
$ws.Range("H286").Value = "B738"
$ws.Range("I286").Value = "(SP-ESD)"
$ws.Range("J286").Value = "11:24 PM"
$ws.Range("L286").Value = "0 hours, 1 minutes"

$ws.Range("A287").Value = 286
$ws.Range("B287").Value = "Saturday, Jan 14"
$ws.Range("C287").Value = "11:52 PM"
$ws.Range("D287").Value = "E41662"
$ws.Range("E287").Value = "Marsa Alam"
$ws.Range("F287").Value = "(RMF)"
$ws.Range("G287").Value = "Enter Air "
$ws.Range("H287").Value = "B738"
$ws.Range("I287").Value = "(SP-ENX)"
$ws.Range("J287").Value = "11:43 PM"
$ws.Range("L287").Value = "0 hours, -9 minutes"

$ws.Range("A288").Value = 287
$ws.Range("B288").Value = "Saturday, Jan 14"
$ws.Range("C288").Value = "11:55 PM"
$ws.Range("D288").Value = "FR7101"
$ws.Range("E288").Value = "Oslo"
$ws.Range("F288").Value = "(OSL)"
$ws.Range("G288").Value = "Ryanair "
$ws.Range("H288").Value = "B738"
$ws.Range("I288").Value = "(SP-RSN)"
$ws.Range("J288").Value = "11:47 PM"
$ws.Range("L288").Value = "0 hours, -8 minutes"

$ws.Range("A289").Value = 288
$ws.Range("B289").Value = "Sunday, Jan 15"
$ws.Range("C289").Value = "12:12 AM"
$ws.Range("D289").Value = "E44092"
$ws.Range("E289").Value = "Marsa Alam"
$ws.Range("F289").Value = "(RMF)"
$ws.Range("G289").Value = "Enter Air "
$ws.Range("H289").Value = "B738"
$ws.Range("I289").Value = "(SP-ENP)"
$ws.Range("J289").Value = "12:21 AM"
$ws.Range("L289").Value = "0 hours, 9 minutes"

$ws.Range("A290").Value = 289
$ws.Range("B290").Value = "Sunday, Jan 15"
$ws.Range("C290").Value = "1:10 AM"
$ws.Range("D290").Value = "W61030"
$ws.Range("E290").Value = "Funchal"
$ws.Range("F290").Value = "(FNC)"
$ws.Range("G290").Value = "Wizz Air "
$ws.Range("H290").Value = "A21N"
$ws.Range("I290").Value = "(HA-LZJ)"
$ws.Range("J290").Value = "12:12 AM"
$ws.Range("L290").Value = "0 hours, -58 minutes"
